# Fix table borders on the header row: the five header cells currently
# only carry an explicit <w:bottom .../> border. Add explicit "nil"
# (no border / not inherited) markers for top, left and right so the
# cell's border set matches the rest of the table, while keeping the
# existing double bottom border untouched.

$d = $word.ActiveDocument
$table = $d.Tables.Item(1)
$headerRow = $table.Rows.Item(1)

for ($i = 1; $i -le $headerRow.Cells.Count; $i++) {
    $cell = $headerRow.Cells.Item($i)

    # wdBorderTop = -1, wdBorderLeft = -2, wdBorderRight = -4
    $cell.Borders.Item(-1).LineStyle = "nil"
    $cell.Borders.Item(-2).LineStyle = "nil"
    $cell.Borders.Item(-4).LineStyle = "nil"
}

Write-Host "Header row borders updated"
